# Altera todos combo de cadastro de veiculo pra autocomplete
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Row 16 ("BD pra carro / deixar cor livre") is now done -> mark "ok" in column C
$ws.Range("C16").Value = "ok"

# Row 17 ("Salvar chat / notificação de chat") gets a new remark in column D
$ws.Range("D17").Value = "complexo"

# Update the active selection to reflect the last edited cell
$ws.Range("D16").Select()
